$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the material name suffix, keep just the label
$ws.Range("A2").Value = "MATERIAL :"

# Remove the buyer name suffix, keep just the label
$ws.Range("A3").Value = "BUYER: "

# Clear out the negociant tag number, shipment number and date lines entirely
$ws.Range("A4").Value = $null
$ws.Range("A5").Value = $null
$ws.Range("A6").Value = $null

# Remove the fixed weight values, keep just the labels
$ws.Range("C10").Value = "NET WEIGHT: "
$ws.Range("C11").Value = "SAMPLE: "
$ws.Range("C12").Value = "DUST:"

# Update the selected cell to match the saved selection state
$ws.Range("C16").Select()
